# Fruta / hortaliza, semanal
# Insert two new weekly report rows at the top of the data block (rows 265-266),
# pushing the previous rows 265-267 down to 267-269.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows above the current row 265 (format/values of the old
# row 265 shift down, carrying its date-style formatting along for the ride).
$ws.Rows("265:266").Insert()

# New row 265: "Especial" quality entry for Región de O'Higgins.
$ws.Range("A265").Value2 = 7
$ws.Range("B265").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C265").Value2 = "Ñuble"
$ws.Range("D265").Value2 = 45075
$ws.Range("E265").Value2 = 16
$ws.Range("F265").Value2 = "Fruta"
$ws.Range("G265").Value2 = 100101
$ws.Range("H265").Value2 = "Berries"
$ws.Range("I265").Value2 = 100101007
$ws.Range("J265").Value2 = "Kiwi"
$ws.Range("K265").Value2 = "Hayward"
$ws.Range("L265").Value2 = "Especial"
$ws.Range("M265").Value2 = 40
$ws.Range("N265").Value2 = 12000
$ws.Range("O265").Value2 = 12000
$ws.Range("P265").Value2 = 12000
$ws.Range("Q265").Value2 = "`$/bandeja 18 kilos"
$ws.Range("R265").Value2 = "Región de O'Higgins"
$ws.Range("S265").Value2 = 667
$ws.Range("T265").Value2 = 18

# New row 266: "Primera" quality entry for Región de O'Higgins.
$ws.Range("A266").Value2 = 7
$ws.Range("B266").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C266").Value2 = "Ñuble"
$ws.Range("D266").Value2 = 45075
$ws.Range("E266").Value2 = 16
$ws.Range("F266").Value2 = "Fruta"
$ws.Range("G266").Value2 = 100101
$ws.Range("H266").Value2 = "Berries"
$ws.Range("I266").Value2 = 100101007
$ws.Range("J266").Value2 = "Kiwi"
$ws.Range("K266").Value2 = "Hayward"
$ws.Range("L266").Value2 = "Primera"
$ws.Range("M266").Value2 = 30
$ws.Range("N266").Value2 = 10000
$ws.Range("O266").Value2 = 10000
$ws.Range("P266").Value2 = 10000
$ws.Range("Q266").Value2 = "`$/bandeja 18 kilos"
$ws.Range("R266").Value2 = "Región de O'Higgins"
$ws.Range("S266").Value2 = 556
$ws.Range("T266").Value2 = 18
